# Updated legacy GSC export data:
# Drop the two oldest date rows (2025-09-08 and 2025-09-09) from the
# "Chart" sheet's video-indexing table, shifting the remaining rows up.
# The newest remaining row (now 2025-09-10) has not finished processing,
# so its "No video indexed" / "Video indexed" counts are cleared to blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the two oldest rows (row 2 = 2025-09-08, row 3 = 2025-09-09).
# This shifts every subsequent row up by two and shrinks the table from
# 88 data rows to 86 data rows (dimension A1:D89 -> A1:D87).
$ws.Range("A2:A3").EntireRow.Delete()

# The row that is now row 2 (2025-09-10) should have blank values for
# "No video indexed" (B) and "Video indexed" (C); "Impressions" (D)
# remains 0. Use a quote-prefixed empty value so the cell is written as
# a (blank) text value, then strip the quote-prefix formatting it adds.
$ws.Range("B2").Value = "'"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = "'"
$ws.Range("C2").ClearFormats()
